# "Generate Report for Handback"
#
# The localization-status workbook gets a handback pass: the Status
# column flips from "Ready for handoff" to "Handed back: in sync with
# en-US" on every sheet, and the per-language sheets (zh-cn / de-de)
# grow two new populated columns - "Latest Target File" (F) and
# "Latest Handback File" (G), each a hyperlinked filename mirroring
# the existing Source/Handoff-file hyperlinks - plus a refreshed
# "Latest Handback DateTime" (H).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet: both per-language status cells flip ----
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

# ---- zh-cn sheet ----
# Re-use the exact display text already used by the existing
# hyperlinks in this row so the new "target"/"handback" file columns
# mirror the "source"/"handoff" file columns.
$zhMdName = $wsZhCn.Range("A2").Value2
$zhXlfName = $wsZhCn.Range("D2").Value2

$wsZhCn.Range("C2").Value = $newStatus

$zhF2 = $wsZhCn.Range("F2")
$wsZhCn.Hyperlinks.Add($zhF2, "https://github.com/OpenLocalizationTestOrg/olhandback/blob/75ea6a315533fc49c6c6682cba8663ceba1811f1/ol-handback/OpenLocalizationTest/oltest/xinjiang/" + $zhMdName, [Type]::Missing, [Type]::Missing, $zhMdName) | Out-Null

$zhG2 = $wsZhCn.Range("G2")
$wsZhCn.Hyperlinks.Add($zhG2, "https://github.com/OpenLocalizationTestOrg/olhandback/blob/79e04b3b9f8d86d204cdee25c82140696f8e5beb/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/hb/" + $zhXlfName, [Type]::Missing, [Type]::Missing, $zhXlfName) | Out-Null

$wsZhCn.Range("H2").Value = "2016-03-25 07:37:02"

# ---- de-de sheet ----
$deMdName = $wsDeDe.Range("A2").Value2
$deXlfName = $wsDeDe.Range("D2").Value2

$wsDeDe.Range("C2").Value = $newStatus

$deF2 = $wsDeDe.Range("F2")
$wsDeDe.Hyperlinks.Add($deF2, "https://github.com/OpenLocalizationTestOrg/olhandback/blob/75ea6a315533fc49c6c6682cba8663ceba1811f1/ol-handback/OpenLocalizationTest/oltest/xinjiang/" + $deMdName, [Type]::Missing, [Type]::Missing, $deMdName) | Out-Null

$deG2 = $wsDeDe.Range("G2")
$wsDeDe.Hyperlinks.Add($deG2, "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a4a98f02f468399ee067ba7553c5312791f08f25/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/hb/" + $deXlfName, [Type]::Missing, [Type]::Missing, $deXlfName) | Out-Null

$wsDeDe.Range("H2").Value = "2016-03-25 07:37:17"
